$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update existing task rows: progress on hours worked / statuses
# ---------------------------------------------------------------------------

# Row 6 - "Insertion dans la base de données": hours revised down from 10 to 1
$ws.Range("E6").Value = 1

# Row 14 - "Chargement du personnage en combat": hours bumped to 4.5, status done
$ws.Range("E14").Value = 4.5
$ws.Range("F14").Value = "Terminé"

# Rows 15 - task finished
$ws.Range("F15").Value = "Terminé"

# Rows 16-18 - tasks now in progress
$ws.Range("F16").Value = "En cours"
$ws.Range("F17").Value = "En cours"
$ws.Range("F18").Value = "En cours"

# Row 21 - "Chronomètre" task finished
$ws.Range("F21").Value = "Terminé"

# ---------------------------------------------------------------------------
# 2. Insert two new task rows (22 & 23) before the "Total" row, pushing the
#    totals block down by two rows.
# ---------------------------------------------------------------------------
$ws.Rows("22:23").Insert()

$ws.Range("A22").Value = 34
$ws.Range("B22").Value = "Fonctionnalité bouton fuir"
$ws.Range("C22").Value = "Fonctionnalité concernant le bouton qui permet de fuir"
$ws.Range("D22").Value = "Anthony Gauthier"
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = "Terminé"

$ws.Range("A23").Value = 35
$ws.Range("B23").Value = "La taverne"
$ws.Range("C23").Value = "Faire la taverne qui regénère la vie"
$ws.Range("D23").Value = "Anthony Gauthier"
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = "Terminé"

# Re-enter the grand-total array formula on its new row (24) so it keeps
# picking up E2:E20 and recalculates with the updated hours above.
$ws.Range("E24").FormulaArray = "=SUM(E2:E20)"
$excel.Calculate()

# ---------------------------------------------------------------------------
# 3. Highlight the new rows' status column the same way as the rest of the
#    "Statut" column (conditional formatting for "Terminé").
# ---------------------------------------------------------------------------
$newRuleRange = $ws.Range("F21:F23")
$newRule = $newRuleRange.FormatConditions.Add(9, 0, $null, $null, "Terminé")
$newRule.SetFirstPriority()
$newRule.Font.Color = 6840149
$newRule.Interior.Color = 13551615

# ---------------------------------------------------------------------------
# 4. Refresh the view: no frozen/scrolled top-left cell anymore, selection
#    rests on E9.
# ---------------------------------------------------------------------------
$ws.Range("E9").Select()
